$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.42516366666666
$ws.Range("H2").Value = 106.275491
$ws.Range("I2").Value = 0.00832770193000585
$ws.Range("J2").Value = 0.008327701930005852
$ws.Range("M2").Value = 18.42392
$ws.Range("N2").Value = 55.27176
$ws.Range("O2").Value = 0.3903243738016154
$ws.Range("P2").Value = 0.3903243738016154
$ws.Range("Q2").Value = 652.6703813815732
$ws.Range("R2").Value = 5874.033432434159
$ws.Range("S2").Value = 0.003250505041036037
$ws.Range("T2").Value = 0.003250505041036039
$ws.Range("G3").Value = 35.42516366666666
$ws.Range("H3").Value = 106.275491
$ws.Range("I3").Value = 0.00832770193000585
$ws.Range("J3").Value = 0.008327701930005852
$ws.Range("O3").Value = 0.2625687066780312
$ws.Range("P3").Value = 0.2625687066780312
$ws.Range("Q3").Value = 439.0471859528743
$ws.Range("R3").Value = 3951.424673575869
$ws.Range("S3").Value = 0.00218659392536178
$ws.Range("T3").Value = 0.00218659392536178
$ws.Range("G4").Value = 35.42516366666666
$ws.Range("H4").Value = 106.275491
$ws.Range("I4").Value = 0.00832770193000585
$ws.Range("J4").Value = 0.008327701930005852
$ws.Range("M4").Value = 10.76369066666667
$ws.Range("N4").Value = 32.291072
$ws.Range("O4").Value = 0.2280367489253622
$ws.Range("P4").Value = 0.2280367489253622
$ws.Range("Q4").Value = 381.3055035240391
$ws.Range("R4").Value = 3431.749531716352
$ws.Range("S4").Value = 0.001899022074137998
$ws.Range("T4").Value = 0.001899022074137998
$ws.Range("G5").Value = 35.42516366666666
$ws.Range("H5").Value = 106.275491
$ws.Range("I5").Value = 0.00832770193000585
$ws.Range("J5").Value = 0.008327701930005852
$ws.Range("M5").Value = 5.620297999999999
$ws.Range("N5").Value = 16.860894
$ws.Range("O5").Value = 0.1190701705949913
$ws.Range("P5").Value = 0.1190701705949913
$ws.Range("Q5").Value = 199.0999765054393
$ws.Range("R5").Value = 1791.899788548954
$ws.Range("S5").Value = 0.0009915808894700345
$ws.Range("T5").Value = 0.0009915808894700347
$ws.Range("I6").Value = 0.01070182047907406
$ws.Range("J6").Value = 0.01070182047907406
$ws.Range("M6").Value = 18.42392
$ws.Range("N6").Value = 55.27176
$ws.Range("O6").Value = 0.3903243738016154
$ws.Range("P6").Value = 0.3903243738016154
$ws.Range("Q6").Value = 838.73814315896
$ws.Range("R6").Value = 7548.64328843064
$ws.Range("S6").Value = 0.004177181377031884
$ws.Range("T6").Value = 0.004177181377031885
$ws.Range("I7").Value = 0.01070182047907406
$ws.Range("J7").Value = 0.01070182047907406
$ws.Range("O7").Value = 0.2625687066780312
$ws.Range("P7").Value = 0.2625687066780312
$ws.Range("S7").Value = 0.002809963162290943
$ws.Range("T7").Value = 0.002809963162290943
$ws.Range("I8").Value = 0.01070182047907406
$ws.Range("J8").Value = 0.01070182047907406
$ws.Range("M8").Value = 10.76369066666667
$ws.Range("N8").Value = 32.291072
$ws.Range("O8").Value = 0.2280367489253622
$ws.Range("P8").Value = 0.2280367489253622
$ws.Range("Q8").Value = 490.0106993135787
$ws.Range("R8").Value = 4410.096293822208
$ws.Range("S8").Value = 0.00244040834963091
$ws.Range("T8").Value = 0.00244040834963091
$ws.Range("I9").Value = 0.01070182047907406
$ws.Range("J9").Value = 0.01070182047907406
$ws.Range("M9").Value = 5.620297999999999
$ws.Range("N9").Value = 16.860894
$ws.Range("O9").Value = 0.1190701705949913
$ws.Range("P9").Value = 0.1190701705949913
$ws.Range("Q9").Value = 255.860767335074
$ws.Range("R9").Value = 2302.746906015666
$ws.Range("S9").Value = 0.001274267590120319
$ws.Range("T9").Value = 0.001274267590120319
$ws.Range("G10").Value = 51.06824600000001
$ws.Range("H10").Value = 153.204738
$ws.Range("I10").Value = 0.01200505761322374
$ws.Range("J10").Value = 0.01200505761322374
$ws.Range("M10").Value = 18.42392
$ws.Range("N10").Value = 55.27176
$ws.Range("O10").Value = 0.3903243738016154
$ws.Range("P10").Value = 0.3903243738016154
$ws.Range("Q10").Value = 940.8772788443201
$ws.Range("R10").Value = 8467.89550959888
$ws.Range("S10").Value = 0.004685866595333871
$ws.Range("T10").Value = 0.004685866595333873
$ws.Range("G11").Value = 51.06824600000001
$ws.Range("H11").Value = 153.204738
$ws.Range("I11").Value = 0.01200505761322374
$ws.Range("J11").Value = 0.01200505761322374
$ws.Range("O11").Value = 0.2625687066780312
$ws.Range("P11").Value = 0.2625687066780312
$ws.Range("Q11").Value = 632.9221202426381
$ws.Range("R11").Value = 5696.299082183743
$ws.Range("S11").Value = 0.003152152451099409
$ws.Range("T11").Value = 0.003152152451099409
$ws.Range("G12").Value = 51.06824600000001
$ws.Range("H12").Value = 153.204738
$ws.Range("I12").Value = 0.01200505761322374
$ws.Range("J12").Value = 0.01200505761322374
$ws.Range("M12").Value = 10.76369066666667
$ws.Range("N12").Value = 32.291072
$ws.Range("O12").Value = 0.2280367489253622
$ws.Range("P12").Value = 0.2280367489253622
$ws.Range("Q12").Value = 549.6828028332375
$ws.Range("R12").Value = 4947.145225499136
$ws.Range("S12").Value = 0.00273759430878121
$ws.Range("T12").Value = 0.00273759430878121
$ws.Range("G13").Value = 51.06824600000001
$ws.Range("H13").Value = 153.204738
$ws.Range("I13").Value = 0.01200505761322374
$ws.Range("J13").Value = 0.01200505761322374
$ws.Range("M13").Value = 5.620297999999999
$ws.Range("N13").Value = 16.860894
$ws.Range("O13").Value = 0.1190701705949913
$ws.Range("P13").Value = 0.1190701705949913
$ws.Range("Q13").Value = 287.018760857308
$ws.Range("R13").Value = 2583.168847715772
$ws.Range("S13").Value = 0.001429444258009249
$ws.Range("T13").Value = 0.001429444258009249
$ws.Range("G14").Value = 4121.876464666667
$ws.Range("H14").Value = 12365.629394
$ws.Range("I14").Value = 0.9689654199776964
$ws.Range("J14").Value = 0.9689654199776964
$ws.Range("M14").Value = 18.42392
$ws.Range("N14").Value = 55.27176
$ws.Range("O14").Value = 0.3903243738016154
$ws.Range("P14").Value = 0.3903243738016154
$ws.Range("Q14").Value = 75941.12223490149
$ws.Range("R14").Value = 683470.1001141134
$ws.Range("S14").Value = 0.3782108207882136
$ws.Range("T14").Value = 0.3782108207882137
$ws.Range("G15").Value = 4121.876464666667
$ws.Range("H15").Value = 12365.629394
$ws.Range("I15").Value = 0.9689654199776964
$ws.Range("J15").Value = 0.9689654199776964
$ws.Range("O15").Value = 0.2625687066780312
$ws.Range("P15").Value = 0.2625687066780312
$ws.Range("Q15").Value = 51085.10661194543
$ws.Range("R15").Value = 459765.9595075088
$ws.Range("S15").Value = 0.254419997139279
$ws.Range("T15").Value = 0.254419997139279
$ws.Range("G16").Value = 4121.876464666667
$ws.Range("H16").Value = 12365.629394
$ws.Range("I16").Value = 0.9689654199776964
$ws.Range("J16").Value = 0.9689654199776964
$ws.Range("M16").Value = 10.76369066666667
$ws.Range("N16").Value = 32.291072
$ws.Range("O16").Value = 0.2280367489253622
$ws.Range("P16").Value = 0.2280367489253622
$ws.Range("Q16").Value = 44366.6032318856
$ws.Range("R16").Value = 399299.4290869703
$ws.Range("S16").Value = 0.2209597241928121
$ws.Range("T16").Value = 0.2209597241928121
$ws.Range("G17").Value = 4121.876464666667
$ws.Range("H17").Value = 12365.629394
$ws.Range("I17").Value = 0.9689654199776964
$ws.Range("J17").Value = 0.9689654199776964
$ws.Range("O17").Value = 0.2625687066780312
$ws.Range("P17").Value = 0.2625687066780312
$ws.Range("Q17").Value = 23166.17405061314
$ws.Range("R17").Value = 208495.5664555182
$ws.Range("S17").Value = 0.1153748778573917
$ws.Range("T17").Value = 0.1153748778573917
